# Apply the "Anonimyzed fedcore" update:
#  - add a top+bottom border to the interior/last cells of the merged
#    header ranges (B1:D1 on both sheets, and E1:G1 on sheet2)
#  - rename the "fedcore" column header to "approach"
#  - drop the stray empty placeholder cell at G5 on the computational
#    comparison sheet

$xlPasteFormats = -4122

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# ---- Build the two new border styles once, on sheet1's C1/D1 ----
# interior cell of the merged block: thin top + thin bottom
$ws1.Range("C1").Style = "Normal"
$ws1.Range("C1").Borders.Item(8).LineStyle = 1    # xlEdgeTop
$ws1.Range("C1").Borders.Item(9).LineStyle = 1    # xlEdgeBottom

# rightmost cell of the merged block = same format + a right edge on top,
# built from a copy of C1 so only one extra border style gets created
$ws1.Range("C1").Copy()
$ws1.Range("D1").PasteSpecial($xlPasteFormats)
$ws1.Range("D1").Borders.Item(10).LineStyle = 1   # xlEdgeRight

# ---- Re-use those exact styles everywhere else they are needed, via
# copy/paste-special so no duplicate/unused style entries are created ----
$ws1.Range("C1").Copy()
$ws2.Range("C1").PasteSpecial($xlPasteFormats)
$ws1.Range("D1").Copy()
$ws2.Range("D1").PasteSpecial($xlPasteFormats)

$ws2.Range("C1").Copy()
$ws2.Range("F1").PasteSpecial($xlPasteFormats)
$ws2.Range("D1").Copy()
$ws2.Range("G1").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# ---- Rename the "fedcore" headers to "approach" ----
$ws1.Range("C2").Value = "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# ---- Remove the stray empty placeholder cell at G5 ----
$ws2.Range("G5").ClearContents()

Write-Host "edit applied"
